$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.02, 1.031612165222644, 1.035000385883583, 1.05286190050676, 1.058456895271841, 1.036942804332752, 1.03674693247731, 1.037798245945326, 1.055609461016022, 1.061189078089102, 1.038219232037962)
    3 = @(1.02, 1.032403668127144, 1.03557408002459, 1.053896420698908, 1.059502416974296, 1.037114613741736, 1.037181187676111, 1.038181871538795, 1.056456343054939, 1.062048049716418, 1.038654103928957)
    4 = @(1.02, 1.032916272654825, 1.035945669075827, 1.054567097668239, 1.060180078957954, 1.037224845135928, 1.037461951801176, 1.038429771178846, 1.057004997429771, 1.062604394332658, 1.038935266771282)
    5 = @(1.02, 1.033131877579069, 1.03610197219917, 1.054849354002281, 1.060465239698899, 1.03727096056853, 1.037579929350078, 1.038533907669481, 1.057235809983088, 1.062838408032877, 1.039053411861836)
    6 = @(1.02, 1.033168084763898, 1.036128221240632, 1.054896763862814, 1.060513135337692, 1.037278690290889, 1.037599735002065, 1.038551387891131, 1.057274573678614, 1.062877707345514, 1.039073245640121)
    7 = @(1.02, 1.032919153163898, 1.035947757265676, 1.05457086799905, 1.060183888223804, 1.037225462220785, 1.037463528443156, 1.03843116297348, 1.057008080937772, 1.062607520739742, 1.038936845652275)
    8 = @(1.02, 1.031879563379807, 1.035194190793469, 1.053211257169191, 1.058809997577962, 1.037001062548586, 1.036893737608318, 1.037927962050159, 1.055895530478439, 1.061479260664259, 1.038366245649095)
    9 = @(1.02, 1.030051183035836, 1.033869220628918, 1.050825248334661, 1.05639780569454, 1.036598463863935, 1.035887995535379, 1.037038764986163, 1.053940210764853, 1.059495248109468, 1.037359075307053)
    10 = @(1.02, 1.028834715474509, 1.032987961035508, 1.049241227628437, 1.054795648821491, 1.036325275773024, 1.035216424661276, 1.03644435590149, 1.052640174500714, 1.058175410034875, 1.036686550725279)
    11 = @(1.02, 1.028308572834153, 1.032606873564321, 1.048556919978159, 1.054103327616705, 1.036205854193634, 1.034925385540147, 1.036186601752237, 1.052078088740697, 1.057604591491191, 1.036395098295215)
    12 = @(1.02, 1.028113231107927, 1.032465398051889, 1.048302976610783, 1.053846383634923, 1.03616132663969, 1.034817245013729, 1.036090805724007, 1.05186943195121, 1.05739266725332, 1.036286804196845)
    13 = @(1.02, 1.028155128424842, 1.032495741503042, 1.048357437495108, 1.053901489246954, 1.036170885591809, 1.034840443118861, 1.036111356753237, 1.051914183784196, 1.057438121036041, 1.036310035245948)
    14 = @(1.02, 1.028292423959164, 1.032595177566985, 1.048535924059503, 1.0540820841559, 1.03620217697743, 1.034916447331175, 1.036178684331409, 1.052060838511608, 1.057587071656743, 1.036386147392961)
    15 = @(1.02, 1.028377028384413, 1.032656453705975, 1.048645927218129, 1.054193383154421, 1.036221434248689, 1.034963271345178, 1.036220159863991, 1.052151214170704, 1.057678858753624, 1.036433037902435)
    16 = @(1.02, 1.028869646515414, 1.033013263311046, 1.049286676328508, 1.054841625981109, 1.036333177647998, 1.035235734911358, 1.036461454470187, 1.052677496025656, 1.058213307796066, 1.03670588839813)
    17 = @(1.02, 1.029178813477369, 1.033237216590038, 1.049689026054351, 1.055248633414312, 1.036402969403839, 1.035406579421274, 1.036612713620499, 1.053007843865707, 1.058548736344132, 1.036876975526847)
    18 = @(1.02, 1.029359202720999, 1.03336789324196, 1.049923862584152, 1.055486171216767, 1.036443568749121, 1.0355062064692, 1.036700904645475, 1.053200611015508, 1.058744451687981, 1.036976744056608)
    19 = @(1.02, 1.029420720477576, 1.03341245877803, 1.050003961648869, 1.05556718871579, 1.036457393569068, 1.035540172670764, 1.036730969387124, 1.053266353316943, 1.058811196649589, 1.037010758494075)
    20 = @(1.02, 1.029145636841381, 1.033213183494805, 1.049645841930366, 1.055204951163913, 1.036395492679467, 1.035388251862012, 1.036596488648349, 1.052972392300208, 1.058512741248231, 1.036858621940348)
    21 = @(1.02, 1.028251991319613, 1.032565893976195, 1.048483357610187, 1.054028897514136, 1.036192967106337, 1.03489406696, 1.036158859537869, 1.052017648839569, 1.057543206554852, 1.036363735239092)
    22 = @(1.02, 1.027690648627259, 1.03215936566267, 1.047753840653772, 1.053290710031975, 1.036064653401034, 1.034583147647925, 1.035883389036893, 1.051418098326142, 1.056934219115764, 1.036052374385931)
    23 = @(1.02, 1.02798817646711, 1.032374830966688, 1.04814043995353, 1.053681918883553, 1.036132767408551, 1.034747990985548, 1.036029450722708, 1.051735861460518, 1.057256997902693, 1.036217451820001)
    24 = @(1.02, 1.029160627757983, 1.033224042866424, 1.049665354505283, 1.055224688868585, 1.036398871426256, 1.03539653337039, 1.036603820125858, 1.052988411086934, 1.058529005680228, 1.036866915209418)
    25 = @(1.02, 1.030523437694394, 1.034211402273196, 1.051440921047433, 1.057020367723815, 1.036703392076119, 1.036148198435192, 1.037268932451497, 1.054445093263997, 1.060007667510912, 1.03761964772483)
}

$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $vals[$i]
    }
}
